$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) holds a date serial of 45243 for rows 2-66.
# Update it to 45244 (one day later) for every row that currently has it.
for ($row = 2; $row -le 66; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
